# logInData.xlsx edit:
#  - clear the placeholder apostrophe ( ' ) entries left in A4/B4/B7/A8
#  - add a new row of real data (AndrewEmad / X123456) in row 9
#  - leave the selection on L9 (clears the stale topLeftCell scroll anchor)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Wipe the quote-prefixed "'" placeholder cells but keep their cell
#    formatting (style index 1 / quotePrefix) intact.
$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B7").ClearContents()

# Row 8 needs the same treatment for A8, but it also carries a stale
# row-level customFormat flag that should disappear once the row is
# rewritten. Capture B8's value, drop the row, then restore B8 and
# re-apply A7's formatting (same style as A8 originally had) onto A8.
$b8Value = $ws.Range("B8").Value2
$ws.Rows("8").Delete()
$ws.Range("B8").Value = $b8Value
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Add the new row of data.
$ws.Range("A9").Value = "AndrewEmad"
$ws.Range("B9").Value = "X123456"

# 3) Move the selection to L9 (also clears the old topLeftCell="A4" pin).
$ws.Range("L9").Select()
